# Macroferia Regional de Talca - Betarraga: add the missing weekly price
# record. A new observation (fecha serial 44551) is inserted as row 184,
# pushing the existing rows 184-228 down to 185-229 and growing the used
# range from A1:R228 to A1:R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 184..228 down to 185..229, leaving a blank row 184 to fill in.
$ws.Rows(184).Insert()

$ws.Range("A184").Value = 5
$ws.Range("B184").Value = "Macroferia Regional de Talca"
$ws.Range("C184").Value = "Maule"
$ws.Range("D184").Value = 44551
$ws.Range("E184").Value = 7
$ws.Range("F184").Value = 100114014
$ws.Range("G184").Value = "Betarraga"
$ws.Range("H184").Value = "Sin especificar"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 5000
$ws.Range("K184").Value = 500
$ws.Range("L184").Value = 500
$ws.Range("M184").Value = 500
$ws.Range("N184").Value = "`$/paquete 5 unidades"
$ws.Range("O184").Value = "Región del Maule"
$ws.Range("P184").Value = 100
$ws.Range("Q184").Value = 5
$ws.Range("R184").Value = "Hortaliza"

# Keep the same date number format as the rest of column D.
$ws.Range("D184").NumberFormat = $ws.Range("D185").NumberFormat
